$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intraday")
$ws.Activate()

# --- New trade-journal rows (143-146) -----------------------------------
# Shared-string insertion order matters (matches how Excel assigned new
# sst indices in the authored file), so write the E-column notes in row
# order (143 -> 281, 144 -> 282, 145 -> 283, 146 -> 284).
$ws.Cells.Item(143, 5).Value = "Missed big holding trade as I went to play cricket in Westbengal"
$ws.Cells.Item(144, 5).Value = "Single trade and logic behind this trade "
$ws.Cells.Item(145, 5).Value = "No trade found as per my setup "
$ws.Cells.Item(146, 5).Value = "early exit due to panic and it gave big return  don't panic in life when you took trade "

# Date cells (column A) reuse the same date style as the rows above them
# (s="1" / numFmtId 14) -- copy the existing format instead of assigning a
# NumberFormat string so no new/duplicate style record gets created.
$dateFormat = $ws.Cells.Item(3, 1)

# Row 143 - 2025-04-01
$ws.Cells.Item(143, 1).Value = 45748
$dateFormat.Copy()
$ws.Cells.Item(143, 1).PasteSpecial(-4122)
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 0

# Row 144 - 2025-04-02
$ws.Cells.Item(144, 1).Value = 45749
$dateFormat.Copy()
$ws.Cells.Item(144, 1).PasteSpecial(-4122)
$ws.Cells.Item(144, 2).Value = 919
$ws.Cells.Item(144, 4).Value = 1

# Row 145 - 2025-04-03
$ws.Cells.Item(145, 1).Value = 45750
$dateFormat.Copy()
$ws.Cells.Item(145, 1).PasteSpecial(-4122)
$ws.Cells.Item(145, 2).Value = " "
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 0

# Row 146 - 2025-04-04
$ws.Cells.Item(146, 1).Value = 45751
$dateFormat.Copy()
$ws.Cells.Item(146, 1).PasteSpecial(-4122)
$ws.Cells.Item(146, 2).Value = 271
$ws.Cells.Item(146, 4).Value = 1

# --- View state: selection / zoom ---------------------------------------
$excel.ActiveWindow.Zoom = 69
$ws.Range("B146").Select()
